# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to match the refreshed data snapshot (commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1248
$ws1.Range("F4").Value  = 17088
$ws1.Range("F6").Value  = 1664
$ws1.Range("F7").Value  = 74
$ws1.Range("F9").Value  = 1027
$ws1.Range("F10").Value = 398
$ws1.Range("F11").Value = 238
$ws1.Range("F12").Value = 134
$ws1.Range("F13").Value = 11828
$ws1.Range("F14").Value = 32
$ws1.Range("F15").Value = 39
$ws1.Range("F16").Value = 2180
$ws1.Range("F17").Value = 4700
$ws1.Range("F18").Value = 491
$ws1.Range("F19").Value = 51
$ws1.Range("F21").Value = 82
$ws1.Range("F23").Value = 346
$ws1.Range("F25").Value = 42

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 1248
$ws4.Range("F4").Value  = 17088
$ws4.Range("F6").Value  = 1664
$ws4.Range("F7").Value  = 74
$ws4.Range("F9").Value  = 1027
$ws4.Range("F10").Value = 398
$ws4.Range("F11").Value = 238
$ws4.Range("F12").Value = 134
$ws4.Range("F15").Value = 11828
$ws4.Range("F16").Value = 32
$ws4.Range("F17").Value = 39
$ws4.Range("F18").Value = 2189
$ws4.Range("F19").Value = 4700
$ws4.Range("F20").Value = 491
$ws4.Range("F21").Value = 51
$ws4.Range("F23").Value = 82
$ws4.Range("F25").Value = 346
$ws4.Range("F27").Value = 42

$wb.Save()
